# Insert a new weekly record row above the current row 312, shifting all
# subsequent rows (312-444) down by one (they become 313-445).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(312).Insert()

# Populate the newly-inserted row 312 with the new weekly record. It carries
# the same market/category/unit metadata as the row that used to occupy this
# slot (now row 313), but with its own date and price figures.
$ws.Range("A312").Value = 5
$ws.Range("B312").Value = "Macroferia Regional de Talca"
$ws.Range("C312").Value = "Maule"
$ws.Range("D312").Value = 44875
$ws.Range("E312").Value = 7
$ws.Range("F312").Value = 100112023
$ws.Range("G312").Value = "Brócoli"
$ws.Range("H312").Value = "Sin especificar"
$ws.Range("I312").Value = "Primera"
$ws.Range("J312").Value = 5000
$ws.Range("K312").Value = 500
$ws.Range("L312").Value = 500
$ws.Range("M312").Value = 500
$ws.Range("N312").Value = "$/unidad"
$ws.Range("O312").Value = "Región del Maule"
$ws.Range("P312").Value = 500
$ws.Range("Q312").Value = 1
$ws.Range("R312").Value = "Hortaliza"
